$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 338; this shifts the existing rows 338-388
# down to 341-391 automatically (including the trailing duplicate block
# that ends up at rows 389-391).
$ws.Rows("338:340").Insert()

# Populate the 3 newly inserted rows (338-340) with the new weekly
# "Hass" / "nueva(o)" price entries (date 2021-09-10 = serial 44449),
# reusing the static descriptor columns shared by every row in this sheet.

# Row 338: Hass - 1a nueva(o)
$ws.Cells.Item(338,1).Value = 8
$ws.Cells.Item(338,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(338,3).Value = "Coquimbo"
$ws.Cells.Item(338,4).Value = 44449
$ws.Cells.Item(338,5).Value = 4
$ws.Cells.Item(338,6).Value = "Fruta"
$ws.Cells.Item(338,7).Value = 100106
$ws.Cells.Item(338,8).Value = "Oleaginosos"
$ws.Cells.Item(338,9).Value = 100106002
$ws.Cells.Item(338,10).Value = "Palta"
$ws.Cells.Item(338,11).Value = "Hass"
$ws.Cells.Item(338,12).Value = "1a nueva(o)"
$ws.Cells.Item(338,13).Value = 400
$ws.Cells.Item(338,14).Value = 2500
$ws.Cells.Item(338,15).Value = 2600
$ws.Cells.Item(338,16).Value = 2550
$ws.Cells.Item(338,17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(338,18).Value = "Provincia de Limar" + [char]0x00ED
$ws.Cells.Item(338,19).Value = 2550
$ws.Cells.Item(338,20).Value = 1

# Row 339: Hass - 2a nueva(o)
$ws.Cells.Item(339,1).Value = 8
$ws.Cells.Item(339,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(339,3).Value = "Coquimbo"
$ws.Cells.Item(339,4).Value = 44449
$ws.Cells.Item(339,5).Value = 4
$ws.Cells.Item(339,6).Value = "Fruta"
$ws.Cells.Item(339,7).Value = 100106
$ws.Cells.Item(339,8).Value = "Oleaginosos"
$ws.Cells.Item(339,9).Value = 100106002
$ws.Cells.Item(339,10).Value = "Palta"
$ws.Cells.Item(339,11).Value = "Hass"
$ws.Cells.Item(339,12).Value = "2a nueva(o)"
$ws.Cells.Item(339,13).Value = 300
$ws.Cells.Item(339,14).Value = 2200
$ws.Cells.Item(339,15).Value = 2300
$ws.Cells.Item(339,16).Value = 2250
$ws.Cells.Item(339,17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(339,18).Value = "Provincia de Limar" + [char]0x00ED
$ws.Cells.Item(339,19).Value = 2250
$ws.Cells.Item(339,20).Value = 1

# Row 340: Hass - 3a nueva (o)
$ws.Cells.Item(340,1).Value = 8
$ws.Cells.Item(340,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(340,3).Value = "Coquimbo"
$ws.Cells.Item(340,4).Value = 44449
$ws.Cells.Item(340,5).Value = 4
$ws.Cells.Item(340,6).Value = "Fruta"
$ws.Cells.Item(340,7).Value = 100106
$ws.Cells.Item(340,8).Value = "Oleaginosos"
$ws.Cells.Item(340,9).Value = 100106002
$ws.Cells.Item(340,10).Value = "Palta"
$ws.Cells.Item(340,11).Value = "Hass"
$ws.Cells.Item(340,12).Value = "3a nueva (o)"
$ws.Cells.Item(340,13).Value = 240
$ws.Cells.Item(340,14).Value = 2000
$ws.Cells.Item(340,15).Value = 2100
$ws.Cells.Item(340,16).Value = 2050
$ws.Cells.Item(340,17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(340,18).Value = "Provincia de Limar" + [char]0x00ED
$ws.Cells.Item(340,19).Value = 2050
$ws.Cells.Item(340,20).Value = 1

# Make sure the date cells carry the same date-time number format as the
# rest of column D.
$ws.Range("D338:D340").NumberFormat = $ws.Range("D341").NumberFormat()
